$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.81798233333333
$ws.Range("H2").Value = 62.453947
$ws.Range("I2").Value = 0.8242653639952813
$ws.Range("J2").Value = 0.8242653639952813
$ws.Range("M2").Value = 7.579746333333333
$ws.Range("N2").Value = 22.739239
$ws.Range("O2").Value = 0.0686314777863378
$ws.Range("P2").Value = 0.0686314777863378
$ws.Range("Q2").Value = 157.7950252584814
$ws.Range("R2").Value = 1420.155227326333
$ws.Range("S2").Value = 0.05657055001908978
$ws.Range("T2").Value = 0.05657055001908978
$ws.Range("G3").Value = 20.81798233333333
$ws.Range("H3").Value = 62.453947
$ws.Range("I3").Value = 0.8242653639952813
$ws.Range("J3").Value = 0.8242653639952813
$ws.Range("O3").Value = 0.0596740760116217
$ws.Range("P3").Value = 0.05967407601162171
$ws.Range("Q3").Value = 137.2004892688592
$ws.Range("R3").Value = 1234.804403419733
$ws.Range("S3").Value = 0.04918727398480145
$ws.Range("T3").Value = 0.04918727398480145
$ws.Range("G4").Value = 20.81798233333333
$ws.Range("H4").Value = 62.453947
$ws.Range("I4").Value = 0.8242653639952813
$ws.Range("J4").Value = 0.8242653639952813
$ws.Range("M4").Value = 96.08192699999999
$ws.Range("N4").Value = 288.245781
$ws.Range("O4").Value = 0.8699822327258658
$ws.Range("P4").Value = 0.8699822327258659
$ws.Range("Q4").Value = 2000.231858838623
$ws.Range("R4").Value = 18002.08672954761
$ws.Range("S4").Value = 0.7170962217272133
$ws.Range("T4").Value = 0.7170962217272134
$ws.Range("G5").Value = 20.81798233333333
$ws.Range("H5").Value = 62.453947
$ws.Range("I5").Value = 0.8242653639952813
$ws.Range("J5").Value = 0.8242653639952813
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.189099
$ws.Range("N5").Value = 0.5672970000000001
$ws.Range("O5").Value = 0.001712213476174646
$ws.Range("P5").Value = 0.001712213476174646
$ws.Range("Q5").Value = 3.936659641251
$ws.Range("R5").Value = 35.429936771259
$ws.Range("S5").Value = 0.00141131826417672
$ws.Range("T5").Value = 0.00141131826417672
$ws.Range("G6").Value = 0.7925996666666667
$ws.Range("I6").Value = 0.03138212158540782
$ws.Range("J6").Value = 0.03138212158540782
$ws.Range("M6").Value = 7.579746333333333
$ws.Range("N6").Value = 22.739239
$ws.Range("O6").Value = 0.0686314777863378
$ws.Range("P6").Value = 0.0686314777863378
$ws.Range("Q6").Value = 6.007704417217889
$ws.Range("R6").Value = 54.069339754961
$ws.Range("S6").Value = 0.002153801380477069
$ws.Range("T6").Value = 0.002153801380477069
$ws.Range("G7").Value = 0.7925996666666667
$ws.Range("I7").Value = 0.03138212158540782
$ws.Range("J7").Value = 0.03138212158540782
$ws.Range("O7").Value = 0.0596740760116217
$ws.Range("P7").Value = 0.05967407601162171
$ws.Range("Q7").Value = 5.223611986973445
$ws.Range("R7").Value = 47.012507882761
$ws.Range("S7").Value = 0.00187269910889358
$ws.Range("T7").Value = 0.001872699108893581
$ws.Range("G8").Value = 0.7925996666666667
$ws.Range("I8").Value = 0.03138212158540782
$ws.Range("J8").Value = 0.03138212158540782
$ws.Range("M8").Value = 96.08192699999999
$ws.Range("N8").Value = 288.245781
$ws.Range("O8").Value = 0.8699822327258658
$ws.Range("P8").Value = 0.8699822327258659
$ws.Range("Q8").Value = 76.15450331289099
$ws.Range("R8").Value = 685.3905298160189
$ws.Range("S8").Value = 0.02730188820454768
$ws.Range("T8").Value = 0.02730188820454768
$ws.Range("G9").Value = 0.7925996666666667
$ws.Range("I9").Value = 0.03138212158540782
$ws.Range("J9").Value = 0.03138212158540782
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.189099
$ws.Range("N9").Value = 0.5672970000000001
$ws.Range("O9").Value = 0.001712213476174646
$ws.Range("P9").Value = 0.001712213476174646
$ws.Range("Q9").Value = 0.149879804367
$ws.Range("R9").Value = 1.348918239303
$ws.Range("S9").Value = 0.00005373289148948651
$ws.Range("T9").Value = 0.00005373289148948651
$ws.Range("G10").Value = 1.536855
$ws.Range("H10").Value = 4.610564999999999
$ws.Range("I10").Value = 0.06085010188305478
$ws.Range("J10").Value = 0.06085010188305479
$ws.Range("M10").Value = 7.579746333333333
$ws.Range("N10").Value = 22.739239
$ws.Range("O10").Value = 0.0686314777863378
$ws.Range("P10").Value = 0.0686314777863378
$ws.Range("Q10").Value = 11.648971051115
$ws.Range("R10").Value = 104.840739460035
$ws.Range("S10").Value = 0.004176232415683266
$ws.Range("T10").Value = 0.004176232415683267
$ws.Range("G11").Value = 1.536855
$ws.Range("H11").Value = 4.610564999999999
$ws.Range("I11").Value = 0.06085010188305478
$ws.Range("J11").Value = 0.06085010188305479
$ws.Range("O11").Value = 0.0596740760116217
$ws.Range("P11").Value = 0.05967407601162171
$ws.Range("Q11").Value = 10.128611628115
$ws.Range("R11").Value = 91.15750465303499
$ws.Range("S11").Value = 0.003631173605084336
$ws.Range("T11").Value = 0.003631173605084336
$ws.Range("G12").Value = 1.536855
$ws.Range("H12").Value = 4.610564999999999
$ws.Range("I12").Value = 0.06085010188305478
$ws.Range("J12").Value = 0.06085010188305479
$ws.Range("M12").Value = 96.08192699999999
$ws.Range("N12").Value = 288.245781
$ws.Range("O12").Value = 0.8699822327258658
$ws.Range("P12").Value = 0.8699822327258659
$ws.Range("Q12").Value = 147.663989919585
$ws.Range("R12").Value = 1328.975909276265
$ws.Range("S12").Value = 0.05293850749781641
$ws.Range("T12").Value = 0.05293850749781642
$ws.Range("G13").Value = 1.536855
$ws.Range("H13").Value = 4.610564999999999
$ws.Range("I13").Value = 0.06085010188305478
$ws.Range("J13").Value = 0.06085010188305479
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.189099
$ws.Range("N13").Value = 0.5672970000000001
$ws.Range("O13").Value = 0.001712213476174646
$ws.Range("P13").Value = 0.001712213476174646
$ws.Range("Q13").Value = 0.290617743645
$ws.Range("R13").Value = 2.615559692805
$ws.Range("S13").Value = 0.0001041883644707666
$ws.Range("T13").Value = 0.0001041883644707666
$ws.Range("G14").Value = 2.108971
$ws.Range("H14").Value = 6.326912999999999
$ws.Range("I14").Value = 0.08350241253625613
$ws.Range("J14").Value = 0.08350241253625615
$ws.Range("M14").Value = 7.579746333333333
$ws.Range("N14").Value = 22.739239
$ws.Range("O14").Value = 0.0686314777863378
$ws.Range("P14").Value = 0.0686314777863378
$ws.Range("Q14").Value = 15.98546520435633
$ws.Range("R14").Value = 143.869186839207
$ws.Range("S14").Value = 0.005730893971087678
$ws.Range("T14").Value = 0.005730893971087679
$ws.Range("G15").Value = 2.108971
$ws.Range("H15").Value = 6.326912999999999
$ws.Range("I15").Value = 0.08350241253625613
$ws.Range("J15").Value = 0.08350241253625615
$ws.Range("O15").Value = 0.0596740760116217
$ws.Range("P15").Value = 0.05967407601162171
$ws.Range("Q15").Value = 13.89913049308967
$ws.Range("R15").Value = 125.092174437807
$ws.Range("S15").Value = 0.004982929312842342
$ws.Range("T15").Value = 0.004982929312842343
$ws.Range("G16").Value = 2.108971
$ws.Range("H16").Value = 6.326912999999999
$ws.Range("I16").Value = 0.08350241253625613
$ws.Range("J16").Value = 0.08350241253625615
$ws.Range("M16").Value = 96.08192699999999
$ws.Range("N16").Value = 288.245781
$ws.Range("O16").Value = 0.8699822327258658
$ws.Range("P16").Value = 0.8699822327258659
$ws.Range("Q16").Value = 202.633997667117
$ws.Range("R16").Value = 1823.705979004053
$ws.Range("S16").Value = 0.07264561529628843
$ws.Range("T16").Value = 0.07264561529628846
$ws.Range("G17").Value = 2.108971
$ws.Range("H17").Value = 6.326912999999999
$ws.Range("I17").Value = 0.08350241253625613
$ws.Range("J17").Value = 0.08350241253625615
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.189099
$ws.Range("N17").Value = 0.5672970000000001
$ws.Range("O17").Value = 0.001712213476174646
$ws.Range("P17").Value = 0.001712213476174646
$ws.Range("Q17").Value = 0.398804307129
$ws.Range("R17").Value = 3.589238764161
$ws.Range("S17").Value = 0.0001429739560376724
$ws.Range("T17").Value = 0.0001429739560376725
